$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1169995834814548, 0.3048912486333797, 3.223369029078222, 0.5333859586016987, 1, 4.178645819794754)
    3 = @(0.1169995834814548, 1.626987699542094, 0.7210945179870265, 13.86384647080068, 1, 16.32892827181126)
    4 = @(0.003078177322033415, 0.04103571897497393, 3.223369029078222, 13.86384647080068, 0, 17.13132939617591)
    5 = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 13.86384647080068, 0, 19.36876847130326)
    6 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 8.656069925401464)
    7 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 1, 3.755628166162433)
    8 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987, 0, 6.15379541431027)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
